$wb = $excel.ActiveWorkbook

$wsM = $wb.Worksheets.Item("M")
$wsA = $wb.Worksheets.Item("A")
$wsQ = $wb.Worksheets.Item("Q")

# Insert a new column H (pushing former H/I/J -> I/J/K) on all three sheets.
$wsM.Columns("H:H").Insert()
$wsA.Columns("H:H").Insert()
$wsQ.Columns("H:H").Insert()

# --- Sheet M (table_name / table_code / series_code touch-ups) ---
$wsM.Range("K4").Value = "umar--mz002--3--M"

$wsM.Range("C5").Value = "sdf"
$wsM.Range("C6").Value = "sdf"

$wsM.Range("C7").Value = "ffd"
$wsM.Range("C8").Value = "ffd"

$wsM.Range("C4").Value = "xcg"

$wsM.Range("J2").Value = "MZ001"
$wsM.Range("J3").Value = "MZ001"

$wsM.Range("J4").Value = "MZ002"

# --- New "interval" column on all three sheets ---
$wsM.Range("H1").Value = "interval"
$wsA.Range("H1").Value = "interval"
$wsQ.Range("H1").Value = "interval"

$wsA.Range("H2:H6").Value = "A"
$wsM.Range("H2:H8").Value = "M"
$wsQ.Range("H2:H5").Value = "Q"

# --- Selection / active cell updates ---
$wsM.Range("H9").Select() | Out-Null
$wsA.Range("I4").Select() | Out-Null
$wsQ.Range("H6").Select() | Out-Null
